$d = $word.ActiveDocument

# Fix the grammar error "it's background" -> "its background" the same way
# Word's grammar-check "Accept" action would: replace the whole affected
# sentence range so the run that used to be split around the flagged word
# (and its now-stale w:proofErr gramStart/gramEnd markers) collapses back
# into a single normal run.
$old = "This is actually going to be pretty easy here. I am going to write the overall narrative of the story, what it" + [char]8217 + "s background is and how the player is going to start out in this game as a whole. This is, of course, by no means the final product and more of a skeleton writing for me to actually collect my thoughts. I am also going to include some names and places as well as factions within and outside of the "
$new = "This is actually going to be pretty easy here. I am going to write the overall narrative of the story, what its background is and how the player is going to start out in this game as a whole. This is, of course, by no means the final product and more of a skeleton writing for me to actually collect my thoughts. I am also going to include some names and places as well as factions within and outside of the "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
